$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new data rows (89 and 90) to the bottom of the table,
# matching the layout of the existing rows (columns A:H numeric,
# column I numeric, column J the string "query").

$ws.Range("A89:H89").Value = 0
$ws.Range("I89").Value = 9.1983
$ws.Range("J89").Value = "query"

$ws.Range("A90:G90").Value = 0
$ws.Range("H90").Value = 0.2
$ws.Range("I90").Value = 9.218299999999999
$ws.Range("J90").Value = "query"
